$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("rating_scales")

$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 1
$ws1.Range("D12").Value = 1
